# Insert a new data row at row 201 (pushing existing rows 201:342 down to
# 202:343, which is exactly what the target diff shows: every former row N
# (N >= 201) reappears unchanged at row N+1, and a brand-new weekly price
# record occupies the newly created row 201).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(201).Insert()

# Populate the freshly inserted row with the new record.
$ws.Cells.Item(201, 1).Value  = 10
$ws.Cells.Item(201, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(201, 3).Value  = "La Araucanía"
$ws.Cells.Item(201, 4).Value  = 44978
$ws.Cells.Item(201, 5).Value  = 9
$ws.Cells.Item(201, 6).Value  = 100112039
$ws.Cells.Item(201, 7).Value  = "Ciboulette"
$ws.Cells.Item(201, 8).Value  = "Sin especificar"
$ws.Cells.Item(201, 9).Value  = "Primera"
$ws.Cells.Item(201, 10).Value = 55
$ws.Cells.Item(201, 11).Value = 6000
$ws.Cells.Item(201, 12).Value = 6000
$ws.Cells.Item(201, 13).Value = 6000
$ws.Cells.Item(201, 14).Value = "`$/docena de atados"
$ws.Cells.Item(201, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(201, 16).Value = 2000
$ws.Cells.Item(201, 17).Value = 3
$ws.Cells.Item(201, 18).Value = "Hortaliza"

# Column D holds dates; make sure the inserted cell keeps the same
# date/time number format used throughout the rest of that column.
$ws.Cells.Item(201, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
